$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as text, forcing text format first for values
# that look numeric (so Excel doesn't auto-convert "0.603" -> 0.603 number).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.122.90"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.922.45"
$ws.Range("E3").Value = "  -4.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "241.16"
$ws.Range("E5").Value = "  -2.47%  "

# Row 6 - XRP
Set-TextValue "D6" "0.603"
$ws.Range("E6").Value = "  -4.39%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.16%  "

# Row 8 - Solana
Set-TextValue "D8" "55.64"
$ws.Range("E8").Value = "  -9.85%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.362"
$ws.Range("E9").Value = "  -5.54%  "

# Row 10 - OKB
Set-TextValue "D10" "54.45"
$ws.Range("E10").Value = "  -4.95%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0824"
$ws.Range("E11").Value = "  +6.17%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.72%  "

# Row 13 - swapped with row 14: now WrappedliquidstakedEther2.0
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.209.72"
$ws.Range("E13").Value = "  -3.97%  "

# Row 14 - swapped with row 13: now Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.806"
$ws.Range("E14").Value = "  -9.39%  "

# Row 15 - Avalanche
Set-TextValue "D15" "20.75"
$ws.Range("E15").Value = "  -7.52%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  -7.36%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -6.73%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "1.938.19"
$ws.Range("E18").Value = "  -3.45%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "36.089.44"
$ws.Range("E19").Value = "  -1.14%  "

# Row 20 - Litecoin
Set-TextValue "D20" "68.59"
$ws.Range("E20").Value = "  -4.69%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  -1.92%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "225.40"
$ws.Range("E22").Value = "  -5.26%  "

# Row 23 - Uniswap
Set-TextValue "D23" "4.92"
$ws.Range("E23").Value = "  -7.13%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.09%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.43"
$ws.Range("E25").Value = "  -2.26%  "

# Row 26 - Toncoin
Set-TextValue "D26" "2.24"
$ws.Range("E26").Value = "  -3.18%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.26"
$ws.Range("E27").Value = "  -5.55%  "

# Row 28 - Monero
Set-TextValue "D28" "161.84"
$ws.Range("E28").Value = "  +1.16%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.04"
$ws.Range("E29").Value = "  -5.44%  "

# Row 30 - Kaspa
Set-TextValue "D30" "0.120"
$ws.Range("E30").Value = "  -12.86%  "

# Row 31 - Stellar: unchanged

# Row 32 - ImmutableX
Set-TextValue "D32" "1.11"
$ws.Range("E32").Value = "  -4.73%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.61"
$ws.Range("E33").Value = "  -7.62%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0618"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue "D35" "4.24"
$ws.Range("E35").Value = "  -4.48%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  -0.04%  "

# Row 37 - THORChain
Set-TextValue "D37" "5.92"
$ws.Range("E37").Value = "  -9.01%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -2.83%  "

# Row 39 - LidoDAOToken
Set-TextValue "D39" "2.11"
$ws.Range("E39").Value = "  -10.09%  "

# Row 40 - RenderToken
Set-TextValue "D40" "2.86"
$ws.Range("E40").Value = "  -9.69%  "

# Row 41 - Cronos
Set-TextValue "D41" "0.0949"
$ws.Range("E41").Value = "  -5.28%  "

# Row 42 - HuobiToken
Set-TextValue "D42" "2.83"
$ws.Range("E42").Value = "  -3.09%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -8.43%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -3.95%  "

# Row 45 - InjectiveProtocol
Set-TextValue "D45" "15.51"
$ws.Range("E45").Value = "  -6.84%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.333.38"
$ws.Range("E46").Value = "  -1.97%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -9.28%  "

# Row 48 - Aave
Set-TextValue "D48" "86.73"
$ws.Range("E48").Value = "  -6.87%  "

# Row 49 - FraxShare
Set-TextValue "D49" "7.15"
$ws.Range("E49").Value = "  -6.01%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  -3.01%  "

# Row 51 - MultiversX
Set-TextValue "D51" "44.99"
$ws.Range("E51").Value = "  +0.90%  "
